# edit.ps1 - applies the LOQ4230.docx content rotation described by the diff.
#
# The underlying paragraph/run structure (styles, bold labels, line breaks,
# italics, headings) is unchanged by the edit: only the *text* assigned to a
# fixed sequence of ten specific runs is cyclically rotated to new slots.
# Because several destination texts equal other slots' original source texts,
# a naive one-pass Find/Replace would corrupt later matches. We therefore use a
# two-phase swap through unique placeholder tokens that cannot collide with any
# existing document text.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Phase 1: move each original run text out of the way into a unique placeholder.
Replace-Text "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia de Produção nos diversos sistemas de produção da indústria. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional." "@@LOQ4230_SLOT0@@"
Replace-Text "To provide an opportunity to apply the fundamental knowledge of Industrial Engineering in the various production systems of the industry. Complementation of general curricular training. Psychological and social adaptation of the student to his future professional activity" "@@LOQ4230_SLOT1@@"
Replace-Text "5840560 - Marco Antonio Carvalho Pereira" "@@LOQ4230_SLOT2@@"
Replace-Text "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais." "@@LOQ4230_SLOT3@@"
Replace-Text "Specific Work Plan. Realization of the Internship. Final and / or partial report." "@@LOQ4230_SLOT4@@"
Replace-Text "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia de Produção. Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio." "@@LOQ4230_SLOT5@@"
Replace-Text "Supervisão das atividades desenvolvidas pelo aluno durante o estágio." "@@LOQ4230_SLOT6@@"
Replace-Text "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio." "@@LOQ4230_SLOT7@@"
Replace-Text "Não será oferecida recuperação." "@@LOQ4230_SLOT8@@"
Replace-Text "A ser definida com o orientador em função das atividades desenvolvidas no estágio." "@@LOQ4230_SLOT9@@"

# Phase 2: drop in the final text for every slot (placeholders are unique, so
# this cannot match anything except the run we just tagged in phase 1).
Replace-Text "@@LOQ4230_SLOT0@@" "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."
Replace-Text "@@LOQ4230_SLOT1@@" "Specific Work Plan. Realization of the Internship. Final and / or partial report."
Replace-Text "@@LOQ4230_SLOT2@@" "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia de Produção nos diversos sistemas de produção da indústria. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."
Replace-Text "@@LOQ4230_SLOT3@@" "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia de Produção. Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."
Replace-Text "@@LOQ4230_SLOT4@@" "To provide an opportunity to apply the fundamental knowledge of Industrial Engineering in the various production systems of the industry. Complementation of general curricular training. Psychological and social adaptation of the student to his future professional activity"
Replace-Text "@@LOQ4230_SLOT5@@" "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
Replace-Text "@@LOQ4230_SLOT6@@" "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio."
Replace-Text "@@LOQ4230_SLOT7@@" "Não será oferecida recuperação."
Replace-Text "@@LOQ4230_SLOT8@@" "A ser definida com o orientador em função das atividades desenvolvidas no estágio."
Replace-Text "@@LOQ4230_SLOT9@@" "5840560 - Marco Antonio Carvalho Pereira"

